$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 4 new rows before existing row 2 (new golf rounds), shifting old data down.
$ws.Rows("2:5").Insert()

# Copy the date-number-format style from an existing date cell so the new
# date cells render the same way (reuses the existing style instead of
# creating a brand new custom number format).
$ws.Range("A19").Copy()
$ws.Range("A2:A5").PasteSpecial(-4122)
$ws.Range("A20").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# --- New row 2: Masterton, Russell, Full-18, Solo ---
$ws.Range("A2").Value = 45956
$ws.Range("B2").Value = "Masterton Golf Course"
$ws.Range("C2").Value = "Russell"
$ws.Range("D2").Value = "Full-18"
$ws.Range("E2").Value = "Solo"
$ws.Range("F2").Value = 115
$ws.Range("G2").Value = 71
$ws.Range("H2").Formula = "=SUM(F2-G2)"
$ws.Range("I2").Value = "I played better at Mahunga, first time playing Masterton back 9, not bad. "

# --- New row 3: Masterton, Hayden, Full-18, Solo ---
$ws.Range("A3").Value = 45956
$ws.Range("B3").Value = "Masterton Golf Course"
$ws.Range("C3").Value = "Hayden"
$ws.Range("D3").Value = "Full-18"
$ws.Range("E3").Value = "Solo"
$ws.Range("F3").Value = 129
$ws.Range("G3").Value = 71
$ws.Range("H3").Formula = "=SUM(F3-G3)"

# --- New row 4: Mahunga, Russell, Full-18, Solo ---
$ws.Range("A4").Value = 45955
$ws.Range("B4").Value = "Mahunga"
$ws.Range("C4").Value = "Russell"
$ws.Range("D4").Value = "Full-18"
$ws.Range("E4").Value = "Solo"
$ws.Range("F4").Value = 104
$ws.Range("G4").Value = 72
$ws.Range("H4").Formula = "=SUM(F4-G4)"
$ws.Range("I4").Value = "Is this my best score from 18? Chips from rough needed work + approach pitches"

# --- New row 5: Mahunga, Hayden, Full-18, Solo ---
$ws.Range("A5").Value = 45955
$ws.Range("B5").Value = "Mahunga"
$ws.Range("C5").Value = "Hayden"
$ws.Range("D5").Value = "Full-18"
$ws.Range("E5").Value = "Solo"
$ws.Range("F5").Value = 121
$ws.Range("G5").Value = 72
$ws.Range("H5").Formula = "=SUM(F5-G5)"
$ws.Range("I5").Value = "Not bad from Hayden, beats his last one"

# --- New row 20 (appended at the bottom): Masterton, Olivia, Front-9, Solo ---
$ws.Range("A20").Value = 45956
$ws.Range("B20").Value = "Masterton Golf Course"
$ws.Range("C20").Value = "Olivia"
$ws.Range("D20").Value = "Front-9"
$ws.Range("E20").Value = "Solo"
$ws.Range("F20").Value = 82
$ws.Range("G20").Value = 36
$ws.Range("H20").Formula = "=SUM(F20-G20)"
$ws.Range("I20").Value = "Was great having Olivia out"

# Refresh the AutoFilter so its range covers the full (now-bigger) table.
$ws.AutoFilterMode = $false
$ws.Range("A1:I20").AutoFilter()

# Keep the _FilterDatabase defined name (driven by the AutoFilter above) in
# sync with the new extent.
$wb.Names.Item("_xlnm._FilterDatabase").RefersTo = "=Summary!`$A`$1:`$I`$20"

# Match the saved selection/active cell.
$ws.Range("G18").Select()
